$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Neutrino Enigma Unraveled: Unveiling the Ghost Particle's Secrets" "Exploring the Realm of Genetics: Unraveling the Secrets of Life"

# Author byline
Replace-Text " Enrico Fermi" " Emily Watson"

# Email address (split across three runs: user / '.' / domain-tld)
Replace-Text "enrico" "emily"
Replace-Text "fermi@physicsinstitute" "watson@schoolofbiology"
Replace-Text "org" "edu"

# Body paragraph, segment before the first double line-break
Replace-Text "In the vast expanse of the universe, there lies a realm of enigmatic particles, elusive and ghostly in their existence--the elusive neutrinos. These subatomic entities, devoid of electric charge and possessing negligible mass, have captivated the scientific community for decades, shrouding them in an aura of mystery and intrigue. Yet, amidst the complexities of their nature and interactions, a breakthrough beckons--a journey into the depths of the neutrino's secrets. Embarking on this voyage of discovery, we delve into the essence of neutrinos, exploring their unique properties, and unraveling the mysteries that have long shrouded them." "In the intricate tapestry of life, genetics holds the key to understanding the symphony of inheritance. From the intricate dance of molecules to the vast canvas of biodiversity, this field unveils the enigmatic secrets of life. It delves into the blueprint of organisms, unravelling the mysteries of how traits and characteristics are passed down through generations."

# Body paragraph, segment between the two double line-breaks
Replace-Text "As neutrinos dance across the cosmos, they possess a remarkable ability to traverse vast distances, passing through matter with ghost-like impunity. Their existence, first postulated by Wolfgang Pauli to explain the energy spectrum of electrons emitted in beta decay, has since been confirmed through meticulous experimentation. Neutrinos come in three distinct flavors, adorned with intriguing names--electron neutrinos, muon neutrinos, and tau neutrinos. Each flavor is associated with its respective charged lepton, exhibiting a fascinating interplay of fundamental particles." "At the heart of genetics lies the study of DNA, the molecule of life. DNA, with its double helix structure, acts as a blueprint for all living organisms. It contains the genetic instructions that determine an organism's traits and characteristics. The study of DNA and its interactions has led to groundbreaking discoveries in fields ranging from medicine to agriculture."

# Body paragraph, segment after the second double line-break
Replace-Text "The elusive nature of neutrinos stems from their extraordinarily weak interactions with other matter. They interact primarily through the electroweak force, one of the four fundamental forces that govern the universe. This feeble interaction allows neutrinos to pass through vast amounts of matter virtually undetected. In fact, trillions of neutrinos from the sun pass through our bodies every second, yet we remain oblivious to their presence. However, this very weakness, which renders neutrinos seemingly intangible, has also hindered our efforts to study and understand them." "Genetics also explores the intricate world of gene expression. Genes, which are segments of DNA, contain the instructions for making proteins. Understanding how genes are expressed and regulated is essential for comprehending a wide range of biological processes, from development and growth to disease and evolution."

# Summary paragraph
Replace-Text "Through extensive research and experimentation, scientists have gained valuable insights into the enigmatic world of neutrinos. Their existence, once merely a hypothesis, is now firmly established, opening avenues for further exploration. The discovery of neutrino flavors and their association with charged leptons has illuminated the intricate tapestry of subatomic interactions. While neutrinos' feeble interactions pose challenges, innovative experimental techniques have emerged, promising to unveil the secrets hidden within these ghostly particles. As we continue to unravel the mysteries of neutrinos, we inch closer to unlocking the fundamental principles that govern our universe and our place within it." "Genetics holds the key to unraveling the mysteries of life, from the intricate dance of DNA to the wonders of biodiversity. Through the study of DNA and gene expression, this field uncovers the secrets of inheritance, variation, and biological processes. Genetics has revolutionized our understanding of life and continues to drive advancements in medicine, agriculture, and biotechnology."

# A new, empty trailing paragraph is appended after the Summary paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
